$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -7
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = -11
$ws.Range("F9").Value = 1
$ws.Range("F13").Value = -2
$ws.Range("F22").Value = 1
$ws.Range("F24").Value = -4
$ws.Range("F25").Value = -3
$ws.Range("F27").Value = -1
